# Add a sample "checked-out" booking row (row 4) plus a few more sample
# rooms (rows 5-7) to the Rooms sheet, so there's data to exercise the
# "checkout date has passed" button.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to smuggle literal text values past the COM layer's
# automatic "looks like a date/number" conversion (Excel would otherwise
# turn "2024-12-06" into a date serial, "044" into 44, etc). We format the
# scratch cell as Text, stuff the literal value in, copy it, and paste
# *values only* into the real destination - PasteSpecial-values leaves the
# destination cell's own number format/style untouched while still landing
# a plain text cell.
$scratch = $ws.Range("Z100")

# NOTE: this COM-interop shell only binds positional parameters, so the
# helper below takes $Address/$Text positionally (no "-Address foo" named
# args).
function Set-LiteralText($Address, $Text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy() | Out-Null
    $ws.Range($Address).PasteSpecial(-4163, $null, $false, $false) | Out-Null
    $scratch.Clear() | Out-Null
}

# --- Row 4: fill in the rest of the existing "333 / Suite" booking ---
$ws.Range("E4").Value = "No"
Set-LiteralText "F4" "2024-12-06"
Set-LiteralText "G4" "2024-12-11"
$ws.Range("H4").Value = "Carlos Alcoba"
Set-LiteralText "I4" "263262662"
$ws.Range("J4").Value = "Pool, Dinner, Room Service"

# --- Row 5: new room ---
Set-LiteralText "A5" "044"
$ws.Range("B5").Value = "Single Room"
Set-LiteralText "C5" "170"
$ws.Range("D5").Value = "City"
$ws.Range("E5").Value = "Yes"

# --- Row 6: new room ---
Set-LiteralText "A6" "028"
$ws.Range("B6").Value = "Suite"
Set-LiteralText "C6" "340"
$ws.Range("D6").Value = "Ocean"
$ws.Range("E6").Value = "Yes"

# --- Row 7: new room ---
Set-LiteralText "A7" "050"
$ws.Range("B7").Value = "Double Room"
Set-LiteralText "C7" "200"
$ws.Range("D7").Value = "Monument"
$ws.Range("E7").Value = "Yes"

# Match the author's final cursor position/selection.
$ws.Range("J7").Select() | Out-Null
